# Update the dSF column (F) values for a handful of rows as part of a
# repull of the data / mean calculation fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -4
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = -6
$ws.Range("F10").Value = 6
